# allure report and code refactor
# - Appends new "used-up" OTP test records to Sheet1 (rows 8-30)
# - Marks the corresponding mobile numbers as "used" on the "Test Data"
#   sheet (column B, rows 33-59)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Test Data")

$sheet1Data = @(
    @{Row=8; Mobile="7980000033"; User="test141236"; Email="test141236@gmail.com"; Pass="SoftSuave107202"},
    @{Row=9; Mobile="7980000037"; User="test84456"; Email="test84456@gmail.com"; Pass="SoftSuave126242"},
    @{Row=10; Mobile="7980000038"; User="test84456"; Email="test84456@gmail.com"; Pass="SoftSuave126242"},
    @{Row=11; Mobile="7980000039"; User="test84456"; Email="test84456@gmail.com"; Pass="SoftSuave126242"},
    @{Row=12; Mobile="7980000040"; User="test115430"; Email="test115430@gmail.com"; Pass="SoftSuave18632"},
    @{Row=13; Mobile="7980000041"; User="test115430"; Email="test115430@gmail.com"; Pass="SoftSuave18632"},
    @{Row=14; Mobile="7980000042"; User="test115430"; Email="test115430@gmail.com"; Pass="SoftSuave18632"},
    @{Row=15; Mobile="7980000043"; User="test77299"; Email="test77299@gmail.com"; Pass="SoftSuave140913"},
    @{Row=16; Mobile="7980000044"; User="test77299"; Email="test77299@gmail.com"; Pass="SoftSuave140913"},
    @{Row=17; Mobile="7980000045"; User="test77299"; Email="test77299@gmail.com"; Pass="SoftSuave140913"},
    @{Row=18; Mobile="7980000046"; User="test88162"; Email="test88162@gmail.com"; Pass="SoftSuave16439"},
    @{Row=19; Mobile="7980000047"; User="test121329"; Email="test121329@gmail.com"; Pass="SoftSuave164883"},
    @{Row=20; Mobile="7980000048"; User="test121329"; Email="test121329@gmail.com"; Pass="SoftSuave164883"},
    @{Row=21; Mobile="7980000049"; User="test121329"; Email="test121329@gmail.com"; Pass="SoftSuave164883"},
    @{Row=22; Mobile="7980000050"; User="test127194"; Email="test127194@gmail.com"; Pass="SoftSuave107168"},
    @{Row=23; Mobile="7980000051"; User="test127194"; Email="test127194@gmail.com"; Pass="SoftSuave107168"},
    @{Row=24; Mobile="7980000052"; User="test127194"; Email="test127194@gmail.com"; Pass="SoftSuave107168"},
    @{Row=25; Mobile="7980000053"; User="test59551"; Email="test59551@gmail.com"; Pass="SoftSuave113390"},
    @{Row=26; Mobile="7980000054"; User="test59551"; Email="test59551@gmail.com"; Pass="SoftSuave113390"},
    @{Row=27; Mobile="7980000055"; User="test59551"; Email="test59551@gmail.com"; Pass="SoftSuave113390"},
    @{Row=28; Mobile="7980000056"; User="test116858"; Email="test116858@gmail.com"; Pass="SoftSuave156281"},
    @{Row=29; Mobile="7980000057"; User="test116858"; Email="test116858@gmail.com"; Pass="SoftSuave156281"},
    @{Row=30; Mobile="7980000058"; User="test116858"; Email="test116858@gmail.com"; Pass="SoftSuave156281"}
)


foreach ($item in $sheet1Data) {
    $r = $item.Row
    # The mobile number looks numeric, so force text storage (matching the
    # other mobile-number cells in the sheet) without leaving a stray
    # NumberFormat applied to the cell afterwards.
    $ws1.Range("A$r").NumberFormat = "@"
    $ws1.Range("A$r").Value = $item.Mobile
    $ws1.Range("A$r").ClearFormats()

    $ws1.Range("B$r").Value = $item.User
    $ws1.Range("C$r").Value = $item.Email
    $ws1.Range("D$r").Value = $item.Pass
}

for ($r = 33; $r -le 59; $r++) {
    $ws2.Range("B$r").Value = "used"
}

Write-Host "Added $($sheet1Data.Count) rows to Sheet1 and marked B33:B59 as used on Test Data"
